$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1644
$ws.Range("I34").Value = 1644
$ws.Range("K34").Value = 1644
$ws.Range("M34").Value = -1441
$ws.Range("H36").Value = 1644
$ws.Range("I36").Value = 1644
$ws.Range("K36").Value = 1644
$ws.Range("M36").Value = -929
$ws.Range("H53").Value = 433.43478
$ws.Range("I53").Value = 164
$ws.Range("K53").Value = 164
$ws.Range("M53").Value = 473
$ws.Range("H69").Value = 6250
$ws.Range("H72").Value = 6250
$ws.Range("H74").Value = 4831.909
$ws.Range("J74").Value = 5975
$ws.Range("L74").Value = 5975
$ws.Range("N74").Value = -7847
$ws.Range("H76").Value = 5874.875
$ws.Range("I76").Value = 4333.1113
$ws.Range("K76").Value = 4333.1113
$ws.Range("M76").Value = -4018.1113
$ws.Range("H77").Value = 4831.909
$ws.Range("J77").Value = 5975
$ws.Range("L77").Value = 29875
$ws.Range("N77").Value = -39235
$ws.Range("H79").Value = 5874.875
$ws.Range("I79").Value = 4333.1113
$ws.Range("K79").Value = 4333.1113
$ws.Range("M79").Value = -3241.1113
$ws.Range("H101").Value = 659.38464
$ws.Range("I101").Value = 305.375
$ws.Range("K101").Value = 916.125
$ws.Range("M101").Value = 705.875
$ws.Range("H112").Value = 2803.6316
$ws.Range("J112").Value = 3855.8333
$ws.Range("L112").Value = 11567.4999
$ws.Range("N112").Value = -13783.4999
$ws.Range("H132").Value = 4210.871
$ws.Range("I132").Value = 4521.4
$ws.Range("J132").Value = 2917
$ws.Range("K132").Value = 13564.2
$ws.Range("L132").Value = 8751
$ws.Range("M132").Value = -11034.2
$ws.Range("N132").Value = -13811
$ws.Range("H137").Value = 3107.6956
$ws.Range("I137").Value = 1432.1333
$ws.Range("K137").Value = 4296.3999
$ws.Range("M137").Value = -1746.3999
$ws.Range("H138").Value = 6243.143
$ws.Range("J138").Value = 16316.5
$ws.Range("L138").Value = 48949.5
$ws.Range("N138").Value = -59229.5
$ws.Range("H140").Value = 112926.664
$ws.Range("I140").Value = 54000
$ws.Range("J140").Value = 124712
$ws.Range("K140").Value = 54000
$ws.Range("L140").Value = 124712
$ws.Range("M140").Value = -48820
$ws.Range("N140").Value = -135072
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11364967
$ws.Range("I32").Value = 13158941
$ws.Range("K32").Value = 13158941
$ws.Range("M32").Value = -13158654
$ws.Range("H43").Value = 250023000
$ws.Range("I43").Value = 1000000000
$ws.Range("K43").Value = 1000000000
$ws.Range("M43").Value = -999999687
$ws.Range("H97").Value = 2183.3845
$ws.Range("I97").Value = 1865.3334
$ws.Range("K97").Value = 1865.3334
$ws.Range("M97").Value = -1369.3334
$ws.Range("H132").Value = 62518056
$ws.Range("I132").Value = 13968.889
$ws.Range("J132").Value = 142880450
$ws.Range("K132").Value = 41906.667
$ws.Range("L132").Value = 428641350
$ws.Range("M132").Value = -39376.667
$ws.Range("N132").Value = -428646410
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1599.08
$ws.Range("I64").Value = 1729.0769
$ws.Range("J64").Value = 1458.25
$ws.Range("K64").Value = 1729.0769
$ws.Range("L64").Value = 1458.25
$ws.Range("M64").Value = -1504.0769
$ws.Range("N64").Value = -1908.25
$ws.Range("H67").Value = 1599.08
$ws.Range("I67").Value = 1729.0769
$ws.Range("J67").Value = 1458.25
$ws.Range("K67").Value = 1729.0769
$ws.Range("L67").Value = 1458.25
$ws.Range("M67").Value = -949.0769
$ws.Range("N67").Value = -3018.25
$ws.Range("H86").Value = 16927.143
$ws.Range("I86").Value = 21501
$ws.Range("J86").Value = 5492.5
$ws.Range("K86").Value = 21501
$ws.Range("L86").Value = 5492.5
$ws.Range("M86").Value = -20378
$ws.Range("N86").Value = -7738.5
$ws.Range("H89").Value = 16927.143
$ws.Range("I89").Value = 21501
$ws.Range("J89").Value = 5492.5
$ws.Range("K89").Value = 107505
$ws.Range("L89").Value = 27462.5
$ws.Range("M89").Value = -101889
$ws.Range("N89").Value = -38694.5
$ws.Range("H105").Value = 35502.668
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23587024
$ws.Range("I31").Value = 1921.85
$ws.Range("J31").Value = 37881024
$ws.Range("K31").Value = 1921.85
$ws.Range("L31").Value = 37881024
$ws.Range("M31").Value = -1626.85
$ws.Range("N31").Value = -37881614
$ws.Range("H34").Value = 23587024
$ws.Range("I34").Value = 1921.85
$ws.Range("J34").Value = 37881024
$ws.Range("K34").Value = 1921.85
$ws.Range("L34").Value = 37881024
$ws.Range("M34").Value = -1719.85
$ws.Range("N34").Value = -37881428
$ws.Range("H38").Value = 17494
$ws.Range("I38").Value = 17494
$ws.Range("K38").Value = 17494
$ws.Range("M38").Value = -17117
$ws.Range("H46").Value = 17494
$ws.Range("I46").Value = 17494
$ws.Range("K46").Value = 17494
$ws.Range("M46").Value = -17283
$ws.Range("H107").Value = 1148.7059
$ws.Range("I107").Value = 947.2
$ws.Range("K107").Value = 947.2
$ws.Range("M107").Value = 972.8
$ws.Range("H132").Value = 1791.303
$ws.Range("I132").Value = 1816.0312
$ws.Range("K132").Value = 5448.0936
$ws.Range("M132").Value = -2918.0936
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 587.1111
$ws.Range("I8").Value = 587.1111
$ws.Range("K8").Value = 1761.3333
$ws.Range("M8").Value = -1622.3333
$ws.Range("H23").Value = 734.6923
$ws.Range("J23").Value = 902.8889
$ws.Range("L23").Value = 2708.6667
$ws.Range("N23").Value = -3178.6667
$ws.Range("H68").Value = 1317.6875
$ws.Range("I68").Value = 1126.7142
$ws.Range("J68").Value = 1466.2222
$ws.Range("K68").Value = 3380.1426
$ws.Range("L68").Value = 4398.6666
$ws.Range("M68").Value = -2569.1426
$ws.Range("N68").Value = -6020.6666
$ws.Range("H71").Value = 1317.6875
$ws.Range("I71").Value = 1126.7142
$ws.Range("J71").Value = 1466.2222
$ws.Range("K71").Value = 10140.4278
$ws.Range("L71").Value = 13195.9998
$ws.Range("M71").Value = -6084.427799999999
$ws.Range("N71").Value = -21307.9998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 21027.75
$ws.Range("J40").Value = 22055.5
$ws.Range("L40").Value = 22055.5
$ws.Range("N40").Value = -22357.5
$ws.Range("H46").Value = 9000
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 24000
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 24000
$ws.Range("N46").Value = -24312
$ws.Range("H113").Value = 4251.4
$ws.Range("I113").Value = 3235.4
$ws.Range("K113").Value = 3235.4
$ws.Range("M113").Value = -1065.4
$ws.Range("H123").Value = 51474.668
$ws.Range("J123").Value = 49549.332
$ws.Range("L123").Value = 49549.332
$ws.Range("N123").Value = -54449.332
$ws.Range("H132").Value = 17795
$ws.Range("I132").Value = 16740.857
$ws.Range("J132").Value = 18465.818
$ws.Range("K132").Value = 50222.571
$ws.Range("L132").Value = 55397.454
$ws.Range("M132").Value = -47692.571
$ws.Range("N132").Value = -60457.454
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 903.1539
$ws.Range("I113").Value = 845.4
$ws.Range("K113").Value = 2536.2
$ws.Range("M113").Value = -366.1999999999998
$ws.Range("H132").Value = 1780.4445
$ws.Range("I132").Value = 1711.409
$ws.Range("J132").Value = 2084.2
$ws.Range("K132").Value = 5134.227000000001
$ws.Range("L132").Value = 6252.599999999999
$ws.Range("M132").Value = -2604.227000000001
$ws.Range("N132").Value = -11312.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N105").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M46").Value = -1344
